$wb = $excel.ActiveWorkbook

# New record data (row 80) to append to both the master "Kayitlar" sheet
# and the district-filtered "Erdemli" sheet.
$recNo     = "3021"
$recDate   = "2025-09-11"
$recBirim  = "Erdemli"
$recParsel = "1"
$recIs     = "3B"
$recPers   = "EMİNE ALANLI KIRCILI (K.Mühendisi), AYHAN KARADAYI (K.Teknisyeni)"

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $row = 80

    $rng = $ws.Range("A" + $row + ":F" + $row)
    # Force text storage so numeric-looking values ("3021", "1") and the
    # date-looking value ("2025-09-11") are kept as plain text, matching
    # the rest of the column (which is stored as text / numberStoredAsText).
    $rng.NumberFormat = "@"

    $ws.Cells.Item($row, 1).Value = $recNo
    $ws.Cells.Item($row, 2).Value = $recDate
    $ws.Cells.Item($row, 3).Value = $recBirim
    $ws.Cells.Item($row, 4).Value = $recParsel
    $ws.Cells.Item($row, 5).Value = $recIs
    $ws.Cells.Item($row, 6).Value = $recPers

    # Remove the explicit number-format style again so the new cells end
    # up with the default (unstyled) appearance, same as the other rows.
    $rng.Style = "Normal"
}
